$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) { Write-Host "MISS para=$paraIndex find=[$findText]" }
}

# Para 1: hyperlink "English" + language list
Replace-InParagraph 1 "English" "ඉංග්‍රීසි"
Replace-InParagraph 1 "Portuguese / French / Thai / Vietnamese / Spanish" "පොච්චුග් / ප්‍රංශ / තයි / වියට්නාම් / ඉස්පාන්"

# Para 3: plain "English"
Replace-InParagraph 3 "English" "ඉංග්‍රීසි"

# Para 5: "Brief"
Replace-InParagraph 5 "Brief" "සාරාංශය"

# Para 6: email description
Replace-InParagraph 6 "An email sent to partners who have attended the event. This email will include a photo gallery It will be sent via customer.io" "ආඩම්බරයට සහභාගී වූ සංස්ථාපකයින්ට යැවුනු ඊමේල් පණිවුඩයක්. මෙම ඊමේලය පින්තූර ගැලරියක් අඩංගු වන අතර customer.io මගින් යැවිය හැක."

# Para 8: "Target audience"
Replace-InParagraph 8 "Target audience" "ලක්ෂ්‍ය ප්‍රේක්ෂකය"

# Para 9: "Event attendees"
Replace-InParagraph 9 "Event attendees" "සැමරීමේ සහභාගීන්"

# Para 12: Subject line pieces
Replace-InParagraph 12 "Subject: " "විෂය: "
Replace-InParagraph 12 "Thank you for coming to " "ඔබට "
Replace-InParagraph 12 "! " " වෙත පැමිණීම ගැන ස්තුතියි! "

# Para 13: success heading
Replace-InParagraph 13 "You made our event a success! 🎉" "ඔබ අපේ උත්සවය සාර්ථක කළා! 🎉"

# Para 15: "Hi "
Replace-InParagraph 15 "Hi " "ආයුබෝවන් "

# Para 17: remove run, replace " in ", replace trailing sentence
Replace-InParagraph 17 "Thank you for attending " ""
Replace-InParagraph 17 " in " " උණුසුම් පිළිගැනීමක් වේවා "
Replace-InParagraph 17 ". We hope you had a great time, and it was a pleasure getting to know you!" " හි ඔබට ස්තුතියි. ඔබ හොඳින් කාලය ගත කළ බව​ අපි අපේක්ෂා කරන්නෙමු, ඔබව දැන හඳුනා ගැනීමට ලැබීම සතුටක්!"

# Para 27: closing "We hope the event..."
Replace-InParagraph 27 "We hope the event inspired you as much as it did us, and let’s keep growing together!" "උත්සවයෙන් අපි දිරිමත් වූ තරමටම ඔබ ද​ දිරිමත් වන්නට​ ඇතැයි අපි බලාපොරොත්තු වෙමු, ​අපි එකට එක්ව ගොඩ නැඟෙමු!"

# Comment: "choose either one"
$c = $d.Comments.Item(1)
$cr = $c.Range
$cr.Find.Execute("choose either one", $true, $false, $false, $false, $false, $true, 1, $false, "ඒකට ඕනෑවක් තෝරන්න", 2) | Out-Null

Write-Host "Edit complete"
